$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2016.1177
$ws.Range("J51").Value = 2160.3076
$ws.Range("L51").Value = 2160.3076
$ws.Range("N51").Value = -3128.3076

$ws.Range("H132").Value = 11914156
$ws.Range("I132").Value = 16675330
$ws.Range("K132").Value = 50025990
$ws.Range("M132").Value = -50023460

$ws.Range("H137").Value = 1231.0847
$ws.Range("I137").Value = 900.7105
$ws.Range("K137").Value = 2702.1315
$ws.Range("M137").Value = -152.1315

$ws.Range("H138").Value = 1210.3562
$ws.Range("I138").Value = 713.3415
$ws.Range("J138").Value = 1847.1562
$ws.Range("K138").Value = 2140.0245
$ws.Range("L138").Value = 5541.4686
$ws.Range("M138").Value = 2999.9755
$ws.Range("N138").Value = -15821.4686

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4456.058
$ws.Range("I32").Value = 3887.7302
$ws.Range("K32").Value = 3887.7302
$ws.Range("M32").Value = -3600.7302

$ws.Range("H61").Value = 23256900
$ws.Range("I61").Value = 27778710
$ws.Range("K61").Value = 27778710
$ws.Range("M61").Value = -27778498

$ws.Range("H74").Value = 1666.2106
$ws.Range("I74").Value = 1176.2667
$ws.Range("K74").Value = 1176.2667
$ws.Range("M74").Value = -302.2666999999999

$ws.Range("H77").Value = 1666.2106
$ws.Range("I77").Value = 1176.2667
$ws.Range("K77").Value = 5881.3335
$ws.Range("M77").Value = -1513.3335

$ws.Range("H136").Value = 23256900
$ws.Range("I136").Value = 27778710
$ws.Range("K136").Value = 83336130
$ws.Range("M136").Value = -83333580

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2494.3333
$ws.Range("I86").Value = 2669.4285
$ws.Range("J86").Value = 1268.6666
$ws.Range("K86").Value = 2669.4285
$ws.Range("L86").Value = 1268.6666
$ws.Range("M86").Value = -1546.4285
$ws.Range("N86").Value = -3514.6666

$ws.Range("H89").Value = 2494.3333
$ws.Range("I89").Value = 2669.4285
$ws.Range("J89").Value = 1268.6666
$ws.Range("K89").Value = 13347.1425
$ws.Range("L89").Value = 6343.333000000001
$ws.Range("M89").Value = -7731.1425
$ws.Range("N89").Value = -17575.333

$ws.Range("H93").Value = 28000
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 28000
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 28000
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = -31744

$ws.Range("H134").Value = 4390.5
$ws.Range("I134").Value = 1219.8148
$ws.Range("J134").Value = 16620.285
$ws.Range("K134").Value = 3659.4444
$ws.Range("L134").Value = 49860.855
$ws.Range("M134").Value = -1124.4444
$ws.Range("N134").Value = -54930.855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1967.3334
$ws.Range("I31").Value = 1924.72
$ws.Range("J31").Value = 2500
$ws.Range("K31").Value = 1924.72
$ws.Range("L31").Value = 2500
$ws.Range("M31").Value = -1629.72
$ws.Range("N31").Value = -3090

$ws.Range("H34").Value = 1967.3334
$ws.Range("I34").Value = 1924.72
$ws.Range("J34").Value = 2500
$ws.Range("K34").Value = 1924.72
$ws.Range("L34").Value = 2500
$ws.Range("M34").Value = -1722.72
$ws.Range("N34").Value = -2904

$ws.Range("H58").Value = 772.21277
$ws.Range("I58").Value = 733.675
$ws.Range("K58").Value = 733.675
$ws.Range("M58").Value = -530.675

$ws.Range("H132").Value = 8053.8423
$ws.Range("I132").Value = 9800.6
$ws.Range("K132").Value = 29401.8
$ws.Range("M132").Value = -26871.8

$ws.Range("H136").Value = 772.21277
$ws.Range("I136").Value = 733.675
$ws.Range("K136").Value = 2201.025
$ws.Range("M136").Value = 348.9750000000004

$ws.Range("H141").Value = 29321.072
$ws.Range("J141").Value = 29321.072
$ws.Range("L141").Value = 29321.072
$ws.Range("N141").Value = -39681.072

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 9289.565000000001
$ws.Range("I3").Value = 4061.111
$ws.Range("J3").Value = 12650.714
$ws.Range("K3").Value = 12183.333
$ws.Range("L3").Value = 37952.142
$ws.Range("M3").Value = -12071.333
$ws.Range("N3").Value = -38176.142

$ws.Range("H5").Value = 2299.0908
$ws.Range("I5").Value = 2429
$ws.Range("K5").Value = 7287
$ws.Range("M5").Value = -7175

$ws.Range("J25").Value = 2000
$ws.Range("L25").Value = 6000
$ws.Range("N25").Value = -6338

$ws.Range("J30").Value = 2000
$ws.Range("L30").Value = 6000
$ws.Range("N30").Value = -6204

$ws.Range("H68").Value = 782.7727
$ws.Range("I68").Value = 872.5
$ws.Range("J68").Value = 675.1
$ws.Range("K68").Value = 2617.5
$ws.Range("L68").Value = 2025.3
$ws.Range("M68").Value = -1806.5
$ws.Range("N68").Value = -3647.3

$ws.Range("H71").Value = 782.7727
$ws.Range("I71").Value = 872.5
$ws.Range("J71").Value = 675.1
$ws.Range("K71").Value = 7852.5
$ws.Range("L71").Value = 6075.900000000001
$ws.Range("M71").Value = -3796.5
$ws.Range("N71").Value = -14187.9

$ws.Range("H80").Value = 2912.8572
$ws.Range("I80").Value = 1020
$ws.Range("J80").Value = 3964.4443
$ws.Range("K80").Value = 3060
$ws.Range("L80").Value = 11893.3329
$ws.Range("M80").Value = -2124
$ws.Range("N80").Value = -13765.3329

$ws.Range("H83").Value = 2912.8572
$ws.Range("I83").Value = 1020
$ws.Range("J83").Value = 3964.4443
$ws.Range("K83").Value = 9180
$ws.Range("L83").Value = 35679.9987
$ws.Range("M83").Value = -4500
$ws.Range("N83").Value = -45039.9987

$ws.Range("H101").Value = 7616.5835
$ws.Range("J101").Value = 7616.5835
$ws.Range("L101").Value = 22849.7505
$ws.Range("N101").Value = -27717.7505

$ws.Range("H135").Value = 2299.0908
$ws.Range("I135").Value = 2429
$ws.Range("K135").Value = 21861
$ws.Range("M135").Value = -19326

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 526.80646
$ws.Range("J107").Value = 342
$ws.Range("L107").Value = 342
$ws.Range("N107").Value = -4182

$ws.Range("H113").Value = 1397.2
$ws.Range("I113").Value = 1333.6
$ws.Range("J113").Value = 1524.4
$ws.Range("K113").Value = 1333.6
$ws.Range("L113").Value = 1524.4
$ws.Range("M113").Value = 836.4000000000001
$ws.Range("N113").Value = -5864.4

$ws.Range("H126").Value = 2921
$ws.Range("I126").Value = 1850
$ws.Range("J126").Value = 3686
$ws.Range("K126").Value = 5550
$ws.Range("L126").Value = 11058
$ws.Range("M126").Value = -3080
$ws.Range("N126").Value = -15998

$ws.Range("H132").Value = 2170.6365
$ws.Range("I132").Value = 1776.5264
$ws.Range("K132").Value = 5329.5792
$ws.Range("M132").Value = -2799.5792

$ws.Range("H134").Value = 28330
$ws.Range("J134").Value = 28330
$ws.Range("L134").Value = 84990
$ws.Range("N134").Value = -90060

$ws.Range("H136").Value = 26555.445
$ws.Range("J136").Value = 26555.445
$ws.Range("L136").Value = 79666.33499999999
$ws.Range("N136").Value = -84766.33499999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 19028.193
$ws.Range("I132").Value = 1310.9348
$ws.Range("J132").Value = 93118.55
$ws.Range("K132").Value = 3932.8044
$ws.Range("L132").Value = 279355.65
$ws.Range("M132").Value = -1402.8044
$ws.Range("N132").Value = -284415.65

$ws.Range("H135").Value = 36087.875
$ws.Range("J135").Value = 36087.875
$ws.Range("L135").Value = 36087.875
$ws.Range("N135").Value = -46227.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 8931909
$ws.Range("I122").Value = 11367602
$ws.Range("J122").Value = 1034.8334
$ws.Range("K122").Value = 34102806
$ws.Range("L122").Value = 3104.5002
$ws.Range("M122").Value = -34100356
$ws.Range("N122").Value = -8004.5002

$ws.Range("H126").Value = 55557296
$ws.Range("I126").Value = 83334760
$ws.Range("K126").Value = 250004280
$ws.Range("M126").Value = -250001810

$ws.Range("H136").Value = 424
$ws.Range("I136").Value = 356.76666
$ws.Range("J136").Value = 712.1429000000001
$ws.Range("K136").Value = 1070.29998
$ws.Range("L136").Value = 2136.4287
$ws.Range("M136").Value = 1479.70002
$ws.Range("N136").Value = -7236.4287
